# labor-timekeeper simulator re-run for Sean Matthew / 2025-12-29 export:
# full-month coverage means this pay period now has 3 timesheet lines
# (Holiday, Regular, OT) instead of 2 (PTO, Regular), and the employee's
# ID got fixed (emp_fn0y5dge -> emp_emnnysju) when the employee record
# was corrected.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Weekly Timesheet" --------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Make room for the extra OT line by inserting a row under the existing
# two detail rows; everything below (subtotal/category/summary block)
# shifts down by one and keeps its per-row styling.
$ws1.Rows("4:4").Insert()

# Row 2: 2026-01-01 / Prezzano / Holiday
$ws1.Range("B2").Value = "Prezzano"
$ws1.Range("C2").Value = 22.5
$ws1.Range("D2").Value = "Holiday"
$ws1.Range("E2").Value = 88
$ws1.Range("F2").Value = 1980

# Row 3: 2026-01-02 / Vincent / Regular
$ws1.Range("A3").Value = "2026-01-02"
$ws1.Range("B3").Value = "Vincent"
$ws1.Range("C3").Value = 17.5
$ws1.Range("D3").Value = "Regular"
$ws1.Range("E3").Value = 88
$ws1.Range("F3").Value = 1540

# Row 4 (new): 2026-01-02 / Vincent / OT
$ws1.Range("A4").Value = "2026-01-02"
$ws1.Range("B4").Value = "Vincent"
$ws1.Range("C4").Value = 5
$ws1.Range("D4").Value = "OT"
$ws1.Range("E4").Value = 88
$ws1.Range("F4").Value = 660

# Row 6 (was row 5): SUBTOTAL line
$ws1.Range("C6").Value = 45
$ws1.Range("D6").Value = "Reg: 40 / OT: 5"
$ws1.Range("F6").Value = 4180

# Row 9 (was row 8): HOURLY SUBTOTAL
$ws1.Range("F9").Value = 4180

# Row 11 (was row 10): GRAND TOTAL
$ws1.Range("F11").Value = 4180

# ---- Sheet 2: "Jason Schema" (persisted simulator log) ----------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

$ws2.Rows("4:4").Insert()

# Row 2: Holiday line for Prezzano
$ws2.Range("B2").Value = "emp_emnnysju"
$ws2.Range("D2").Value = "Prezzano"
$ws2.Range("E2").Value = 22.5
$ws2.Range("F2").Value = 88
$ws2.Range("G2").Value = 1980
$ws2.Range("H2").Value = "Holiday"
$ws2.Range("I2").Value = ""

# Row 3: Regular line for Vincent
$ws2.Range("B3").Value = "emp_emnnysju"
$ws2.Range("C3").Value = "2026-01-02"
$ws2.Range("D3").Value = "Vincent"
$ws2.Range("E3").Value = 17.5
$ws2.Range("F3").Value = 88
$ws2.Range("G3").Value = 1540
$ws2.Range("H3").Value = "Regular"
$ws2.Range("I3").Value = ""

# Row 4 (new): OT line for Vincent
$ws2.Range("A4").Value = "Sean Matthew"
$ws2.Range("B4").Value = "emp_emnnysju"
$ws2.Range("C4").Value = "2026-01-02"
$ws2.Range("D4").Value = "Vincent"
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 88
$ws2.Range("G4").Value = 660
$ws2.Range("H4").Value = "OT"
$ws2.Range("I4").Value = ""
